# Apply the "Exceptional items" column insertion to the Quarterly sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a new blank column at L, shifting existing columns L:T right to M:U.
$ws.Columns("L:L").Insert()

# Fill in the new column's header labels (row 1 = lowercase style header,
# row 2 = capitalized style header used elsewhere on this sheet).
$ws.Range("L1").Value2 = "Exceptional items"
$ws.Range("L2").Value2 = "Exceptional Items"

# Only quarter "Mar '16" (row 16) has an actual reported exceptional item;
# all the other data rows keep the new column blank.
$ws.Range("L16").Value2 = 18.19
